$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the six previously-numeric U-column cells to the text "NA"
# (these held ad-hoc placeholder values that got replaced when the sheet
# was restored to the original/upstream metadata content).
$ws.Cells.Item(7, 21).Value2 = "NA"
$ws.Cells.Item(16, 21).Value2 = "NA"
$ws.Cells.Item(28, 21).Value2 = "NA"
$ws.Cells.Item(29, 21).Value2 = "NA"
$ws.Cells.Item(34, 21).Value2 = "NA"
$ws.Cells.Item(39, 21).Value2 = "NA"

# Append the restored Time_to_death_or_last_followup_days values (column U)
# for the additional rows 69-113 that exist in the original metadata file.
$ws.Cells.Item(69, 21).Value2 = 1461
$ws.Cells.Item(70, 21).Value2 = 1965
$ws.Cells.Item(71, 21).Value2 = 370
$ws.Cells.Item(72, 21).Value2 = 1205
$ws.Cells.Item(73, 21).Value2 = 695
$ws.Cells.Item(74, 21).Value2 = 1028
$ws.Cells.Item(75, 21).Value2 = 1161
$ws.Cells.Item(76, 21).Value2 = 864
$ws.Cells.Item(77, 21).Value2 = 808
$ws.Cells.Item(78, 21).Value2 = 980
$ws.Cells.Item(79, 21).Value2 = 737
$ws.Cells.Item(80, 21).Value2 = 1494
$ws.Cells.Item(81, 21).Value2 = 993
$ws.Cells.Item(82, 21).Value2 = 895
$ws.Cells.Item(83, 21).Value2 = 2360
$ws.Cells.Item(84, 21).Value2 = 1036
$ws.Cells.Item(85, 21).Value2 = 1992
$ws.Cells.Item(86, 21).Value2 = 1163
$ws.Cells.Item(87, 21).Value2 = 4085
$ws.Cells.Item(88, 21).Value2 = 1393
$ws.Cells.Item(89, 21).Value2 = 453
$ws.Cells.Item(90, 21).Value2 = 833
$ws.Cells.Item(91, 21).Value2 = 2004
$ws.Cells.Item(92, 21).Value2 = 905
$ws.Cells.Item(93, 21).Value2 = 1109
$ws.Cells.Item(94, 21).Value2 = 724
$ws.Cells.Item(95, 21).Value2 = 496
$ws.Cells.Item(96, 21).Value2 = 708
$ws.Cells.Item(97, 21).Value2 = 2037
$ws.Cells.Item(98, 21).Value2 = 1365
$ws.Cells.Item(99, 21).Value2 = 702
$ws.Cells.Item(100, 21).Value2 = 2176
$ws.Cells.Item(101, 21).Value2 = 1584
$ws.Cells.Item(102, 21).Value2 = 2068
$ws.Cells.Item(103, 21).Value2 = 1017
$ws.Cells.Item(104, 21).Value2 = 2146
$ws.Cells.Item(105, 21).Value2 = 682
$ws.Cells.Item(106, 21).Value2 = 2214
$ws.Cells.Item(107, 21).Value2 = 1136
$ws.Cells.Item(108, 21).Value2 = 647
$ws.Cells.Item(109, 21).Value2 = 1993
$ws.Cells.Item(110, 21).Value2 = 188
$ws.Cells.Item(111, 21).Value2 = 1868
$ws.Cells.Item(112, 21).Value2 = 997
$ws.Cells.Item(113, 21).Value2 = 910

# New best-fit column width for column U (21) to accommodate the data.
$ws.Columns.Item(21).ColumnWidth = 33

# Restore the original selection/view state on the sheet.
$ws.Range("U1:U1048576").Select()
